$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rangeToSort = $ws.Range("A2:D18")
$keyColumn = $ws.Range("A2:A18")

$rangeToSort.Sort($keyColumn, 1)
